# testPlanTab_TestData.xlsx update
#
# Commit: "Changed locator and worked on author"
#
# Functional edit captured by this script:
#   - tc002!A2 release-name test value is updated from the old
#     "STG- PulseCodeOnAzureCloud" locator/test value to the new
#     "STG- SPARK Modernization" value (a new shared string is created).
#   - The author/tester then left their selection on tc002!A7 with tc002
#     as the active sheet/tab (previously tc010 was the active tab with
#     its own stale selection) - this is the normal "tabSelected" /
#     active-cell bookkeeping that Excel persists for whichever sheet and
#     cell were last touched/selected before saving.

$wb = $excel.ActiveWorkbook

# --- Update the release/locator value on the "tc002" sheet ---------------
$ws = $wb.Worksheets.Item("tc002")
$ws.Activate()

$ws.Range("A2").Value = "STG- SPARK Modernization"

# Leave the selection where the author left it when they saved the file.
$ws.Range("A7").Select()
